$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the column headers: "<name>_old" -> "<name>_FV2304" and
#    "<name>_new" -> "<name>_FV2310" (row 1, columns A:U).
# ---------------------------------------------------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2304"
    "B1" = "Segmentgruppe_FV2304"
    "C1" = "Segment_FV2304"
    "D1" = "Datenelement_FV2304"
    "E1" = "Segment ID_FV2304"
    "F1" = "Code_FV2304"
    "G1" = "Qualifier_FV2304"
    "H1" = "Beschreibung_FV2304"
    "I1" = "Bedingungsausdruck_FV2304"
    "J1" = "Bedingung_FV2304"
    "L1" = "Segmentname_FV2310"
    "M1" = "Segmentgruppe_FV2310"
    "N1" = "Segment_FV2310"
    "O1" = "Datenelement_FV2310"
    "P1" = "Segment ID_FV2310"
    "Q1" = "Code_FV2310"
    "R1" = "Qualifier_FV2310"
    "S1" = "Beschreibung_FV2310"
    "T1" = "Bedingungsausdruck_FV2310"
    "U1" = "Bedingung_FV2310"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# ---------------------------------------------------------------------------
# 2) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true | Out-Null

# ---------------------------------------------------------------------------
# 3) Turn the used range into a native Excel table ("Table1").
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U63"), $null, 1)
$lo.Name = "Table1"

Write-Output "done"
